$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Combat 1/2 mob names: old generic names -> new character names ---
$ws.Range("M6").Value = "Ella Jarvis"
$ws.Range("M7").Value = "Maggie Ortega"
$ws.Range("M10").Value = "Puck Yaztromo"
$ws.Range("M13").Value = "Ron Doom"
$ws.Range("M11").Value = "Danica Lacroix"
$ws.Range("M12").Value = "Isobel Flowright"
$ws.Range("M4").Value = "Wolf I"
$ws.Range("M5").Value = "Wolf II"

# --- Swap the rank numbers of rows 11 and 12 ---
$ws.Range("L11").Value = 3
$ws.Range("L12").Value = 4

# --- Update "weak" elements for combat 2 rows (U10/U11 -> earth, U12 -> water) ---
$ws.Range("U10").Value = "earth"
$ws.Range("U11").Value = "earth"
$ws.Range("U12").Value = "water"

# Match the cell fill/border formatting to the new grouping: rows 10 & 11 both
# "earth" take on the style already used by row 6 / row 12's old "earth" cell;
# row 12 "water" takes on the style already used by row 7 / row 13's "water" cell.
$ws.Range("U6").Copy()
$ws.Range("U10").PasteSpecial(-4122)
$ws.Range("U6").Copy()
$ws.Range("U11").PasteSpecial(-4122)
$ws.Range("U7").Copy()
$ws.Range("U12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column M a bit wider to fit the longer character names ---
$ws.Columns("M").ColumnWidth = 17.16

# --- Selection moves to J18 ---
$null = $ws.Range("J18").Select()
